# Generate Report for Handoff
# The localization run moved from "In Translation" to "Ready for handoff";
# refresh the Status + timestamp columns on every sheet to reflect the new
# handoff generation, then let Excel re-flow the now-wider Status column.

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$overviewStamp   = "2016-09-06 02:43:14"   # Latest HO Xliff Generate Date (Overview) / Latest Handoff Datetime (de-de)
$zhHandoffStamp  = "2016-09-06 02:43:09"   # Latest Handoff Datetime (zh-cn)

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus   # zh-cn status
$wsOverview.Range("F2").Value = $newStatus   # de-de status
$wsOverview.Range("G2").Value = $overviewStamp

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $zhHandoffStamp

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $overviewStamp

# --- Re-fit the Status columns now that "Ready for handoff" is wider -------
# Excel auto-sized these columns to the new text on handoff generation
# (matches the wider "Ready for handoff" label replacing "In Translation").
$statusColWidth = 16.33
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth  # Overview column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth  # Overview column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth      # zh-cn column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth      # de-de column C (Status)
